$wb = $excel.ActiveWorkbook

# The "建物" (Building) sheet is the 2nd sheet in the workbook.
$ws = $wb.Worksheets.Item(2)

# Column I holds "property_category"; it was incorrectly left as "land"
# (copied from the land-sheet template) for every data row. Correct it
# to "building" for all data rows (2-9).
$ws.Range("I2:I9").Value = "building"
